$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '29.415.07'
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '1.915.17'
$ws.Range('E3').Value = '  +1.04%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '325.04'
$ws.Range('E5').Value = '  +0.32%  '
$ws.Range('E6').Value = '  +0.50%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.4818'
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.4075'
$ws.Range('E8').Value = '  +0.43%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.08231'
$ws.Range('E9').Value = '  +2.53%  '
$ws.Range('E10').Value = '  +1.51%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '23.47'
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '1.915.94'
$ws.Range('E12').Value = '  +0.71%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '6.090'
$ws.Range('E13').Value = '  +2.92%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '7.245'
$ws.Range('E14').Value = '  +2.72%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '91.26'
$ws.Range('E15').Value = '  +2.00%  '
$ws.Range('E16').Value = '  +2.13%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '1.008'
$ws.Range('E17').Value = '  +0.56%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '0.00001039'
$ws.Range('E18').Value = '  +1.51%  '
$ws.Range('E19').Value = '  +0.60%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '1.006'
$ws.Range('E20').Value = '  +0.54%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '29.446.85'
$ws.Range('E21').Value = '  +0.23%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '5.661'
$ws.Range('E22').Value = '  +2.63%  '
$ws.Range('E23').Value = '  +1.13%  '
$ws.Range('E24').Value = '  +1.16%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '2.174.72'
$ws.Range('E25').Value = '  +2.21%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '6.576'
$ws.Range('E26').Value = '  +9.48%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '155.92'
$ws.Range('E27').Value = '  +0.95%  '
$ws.Range('E28').Value = '  +1.43%  '
$ws.Range('E29').Value = '  +1.45%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '120.32'
$ws.Range('E30').Value = '  +2.21%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '1.021'
$ws.Range('E31').Value = '  +0.31%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '0.09578'
$ws.Range('E32').Value = '  +1.21%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '5.688'
$ws.Range('E33').Value = '  +6.29%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '3.548'
$ws.Range('E34').Value = '  +0.49%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '1.372'
$ws.Range('E35').Value = '  -0.86%  '
$ws.Range('E36').Value = '  +1.93%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.06111'
$ws.Range('E37').Value = '  +1.29%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '1.182'
$ws.Range('E38').Value = '  +1.07%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.5995'
$ws.Range('E39').Value = '  +2.41%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '8.069'
$ws.Range('E40').Value = '  +3.16%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '10.79'
$ws.Range('E41').Value = '  +6.92%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.1850'
$ws.Range('E42').Value = '  +0.72%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '2.424'
$ws.Range('E43').Value = '  +0.57%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '1.276'
$ws.Range('E44').Value = '  -0.82%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.07625'
$ws.Range('E45').Value = '  -1.14%  '
$ws.Range('E46').Value = '  +1.79%  '
$ws.Range('E47').Value = '  +1.61%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '1.958'
$ws.Range('E48').Value = '  +2.18%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '117.64'
$ws.Range('E49').Value = '  +4.18%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '2.425'
$ws.Range('E50').Value = '  +3.96%  '
$ws.Range('E51').Value = '  +1.15%  '
